$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet is protected; unprotect to allow cell edits, then re-protect afterwards.
$ws.Unprotect()

# Update the "as of" date in the confidential disclaimer note (A10).
$newNote = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-05-24 for illustrative purposes only and are subject to change."
$ws.Range("A10").Value = $newNote
# Re-fit the row height (setting a multi-line value can mark the row as a
# custom height); AutoFit restores it to the sheet's normal computed height.
$ws.Rows(10).AutoFit()

# Update Weight (D) and Percent Change (E) values for the holdings table.
$ws.Range("D2").Value = 0.4777752607154719
$ws.Range("E2").Value = 0.003125000000000044

$ws.Range("D3").Value = 0.3405302897752631
$ws.Range("E3").Value = 0.003656135745757849

$ws.Range("D4").Value = 0.09639980280965545
$ws.Range("E4").Value = 0.01212553495007129

$ws.Range("D5").Value = 0.05347836523556326
$ws.Range("E5").Value = 0.00114797382619658

$ws.Range("D6").Value = 0.03181628146404622
$ws.Range("E6").Value = 0.01331403762662808

$ws.Range("E7").Value = 0.004391966764954791

# Restore sheet protection as it was before editing.
$ws.Protect("D382")
